$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 9216
$ws.Range("F7").Value = 127
$ws.Range("F8").Value = 1281
$ws.Range("F10").Value = 663
$ws.Range("F13").Value = 159
$ws.Range("F14").Value = 300
$ws.Range("F15").Value = 8
$ws.Range("F16").Value = 58
$ws.Range("F17").Value = 1540
$ws.Range("F19").Value = 572
$ws.Range("F21").Value = 1407
$ws.Range("F22").Value = 93
$ws.Range("F23").Value = 250
$ws.Range("F26").Value = 75
$ws.Range("F28").Value = 329
$ws.Range("F29").Value = 329
$ws.Range("F32").Value = 37
$ws.Range("F34").Value = 218
$ws.Range("F36").Value = 585
$ws.Range("F37").Value = 615
$ws.Range("F42").Value = 103
$ws.Range("F43").Value = 524
$ws.Range("F45").Value = 700
$ws.Range("F46").Value = 241
$ws.Range("F48").Value = 50

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 675
$ws.Range("F21").Value = 9
$ws.Range("F23").Value = 124
$ws.Range("F26").Value = 1039
$ws.Range("F27").Value = 249
$ws.Range("F30").Value = 254

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 331
$ws.Range("F6").Value = 150
$ws.Range("F7").Value = 2208
$ws.Range("F8").Value = 3297
$ws.Range("F9").Value = 37

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 9216
$ws.Range("F6").Value = 331
$ws.Range("F7").Value = 150
$ws.Range("F8").Value = 2208
$ws.Range("F9").Value = 3298
$ws.Range("F10").Value = 127
$ws.Range("F11").Value = 1281
$ws.Range("F12").Value = 663
$ws.Range("F15").Value = 159
$ws.Range("F16").Value = 300
$ws.Range("F17").Value = 1540
$ws.Range("F18").Value = 675
$ws.Range("F20").Value = 572
$ws.Range("F21").Value = 37
$ws.Range("F22").Value = 1407
$ws.Range("F23").Value = 93
$ws.Range("F24").Value = 250
$ws.Range("F27").Value = 75
$ws.Range("F28").Value = 329
$ws.Range("F29").Value = 329
$ws.Range("F33").Value = 9
$ws.Range("F34").Value = 37
$ws.Range("F38").Value = 218
$ws.Range("F39").Value = 1039
$ws.Range("F40").Value = 249
$ws.Range("F41").Value = 585
$ws.Range("F42").Value = 615
$ws.Range("F45").Value = 254
$ws.Range("F47").Value = 103
$ws.Range("F49").Value = 524
$ws.Range("F50").Value = 700
